# UseCaseCourseData.xlsx - "Fixed UseCase Data and Creation"
#
# On the "Survey" sheet, column BC ("AssessmentType") was wrong: row 2 held a
# bare numeric placeholder (1) and rows 3-14 had no AssessmentType at all.
# Every question row is actually a SURVEY-type item, so column BC is filled
# in for rows 2-14 with the "SURVEY" label (already present in the shared
# string table), rendered in a small custom font/color and vertically
# centered - matching the existing BC1 header styling convention used
# elsewhere on this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Survey")

# Fill in the AssessmentType column for every data row.
$ws.Range("BC2:BC14").Value = "SURVEY"

# Build the "SURVEY" cell style once on BC2 ...
$cell = $ws.Range("BC2")
$cell.VerticalAlignment = -4108          # xlVAlignCenter
$cell.Font.Size = 8
$cell.Font.Name = "MesloLGM NF"
$cell.Font.Family = 3                    # Modern (monospace) font family
$cell.Font.Color = 7185097               # RGB(201, 162, 109) = #C9A26D

# ... then stamp the same formatting onto the rest of the column.
$cell.Copy()
$ws.Range("BC3:BC14").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Leave the sheet scrolled/selected the way it was left in the source edit.
$ws.Activate()
$ws.Range("BD1").Select()
